# Weekly update: insert a new price record row at row 212 for
# "Vega Modelo de Temuco - Cilantro", shifting the existing rows 212:246
# down to 213:247 (dimension grows from A1:R246 to A1:R247).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 212; this pushes the
# previous rows 212-246 down to 213-247 automatically.
$ws.Rows.Item(212).Insert()

# Populate the newly inserted row 212 with the new data point.
$ws.Cells.Item(212, 1).Value  = 10
$ws.Cells.Item(212, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(212, 3).Value  = "La Araucanía"
$ws.Cells.Item(212, 4).Value  = 44491
$ws.Cells.Item(212, 5).Value  = 9
$ws.Cells.Item(212, 6).Value  = 100112040
$ws.Cells.Item(212, 7).Value  = "Cilantro"
$ws.Cells.Item(212, 8).Value  = "Sin especificar"
$ws.Cells.Item(212, 9).Value  = "Primera"
$ws.Cells.Item(212, 10).Value = 110
$ws.Cells.Item(212, 11).Value = 6000
$ws.Cells.Item(212, 12).Value = 6000
$ws.Cells.Item(212, 13).Value = 6000
$ws.Cells.Item(212, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(212, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(212, 16).Value = 3000
$ws.Cells.Item(212, 17).Value = 2
$ws.Cells.Item(212, 18).Value = "Hortaliza"

# Make sure the date cell keeps the same date-formatted style as the
# rest of column D (style index "2" in styles.xml).
$ws.Cells.Item(212, 4).NumberFormat = $ws.Cells.Item(213, 4).NumberFormat
